$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data, shifting everything down.
$ws.Rows.Item(1).Insert()

# Populate the new header row: Name | Description | Description_FR
# (written in this order so new shared-string indices land as Name, Description_FR, Description)
$ws.Range("A1").Value = "Name"
$ws.Range("C1").Value = "Description_FR"
$ws.Range("B1").Value = "Description"

# Bold the header row.
$ws.Range("A1:C1").Font.Bold = $true

# Match the resulting selection state.
$ws.Range("B1").Select()
